$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "胜宏科技"
$ws.Range("B2").Value = "工业富联"
$ws.Range("C2").Value = "利欧股份"
$ws.Range("A3").Value = "工业富联"
$ws.Range("B3").Value = "利欧股份"
$ws.Range("C3").Value = "先导智能"
$ws.Range("A4").Value = "国轩高科"
$ws.Range("B4").Value = "三维通信"
$ws.Range("C4").Value = "步步高"
$ws.Range("A5").Value = "三维通信"
$ws.Range("B5").Value = "胜宏科技"
$ws.Range("C5").Value = "春兴精工"
$ws.Range("A6").Value = "利欧股份"
$ws.Range("B6").Value = "春兴精工"
$ws.Range("C6").Value = "国轩高科"
$ws.Range("A7").Value = "寒武纪-U"
$ws.Range("B7").Value = "中国联通"
$ws.Range("C7").Value = "卧龙电驱"
$ws.Range("A8").Value = "新易盛"
$ws.Range("B8").Value = "国轩高科"
$ws.Range("C8").Value = "幸福蓝海"
$ws.Range("A9").Value = "三江购物"
$ws.Range("B9").Value = "首开股份"
$ws.Range("C9").Value = "吉视传媒"
$ws.Range("A10").Value = "中国联通"
$ws.Range("B10").Value = "吉视传媒"
$ws.Range("C10").Value = "胜宏科技"
$ws.Range("A11").Value = "吉视传媒"
$ws.Range("B11").Value = "万马股份"
$ws.Range("C11").Value = "阳光电源"
$ws.Range("A12").Value = "先导智能"
$ws.Range("B12").Value = "先导智能"
$ws.Range("C12").Value = "天际股份"
$ws.Range("A13").Value = "春兴精工"
$ws.Range("B13").Value = "寒武纪-U"
$ws.Range("C13").Value = "万通发展"
$ws.Range("A14").Value = "中际旭创"
$ws.Range("B14").Value = "步步高"
$ws.Range("C14").Value = "元隆雅图"
$ws.Range("A15").Value = "卧龙电驱"
$ws.Range("B15").Value = "岩山科技"
$ws.Range("C15").Value = "岩山科技"
$ws.Range("A16").Value = "立讯精密"
$ws.Range("B16").Value = "天际股份"
$ws.Range("C16").Value = "三维通信"
$ws.Range("A17").Value = "步步高"
$ws.Range("B17").Value = "三江购物"
$ws.Range("C17").Value = "巨人网络"
$ws.Range("A18").Value = "岩山科技"
$ws.Range("B18").Value = "中际旭创"
$ws.Range("C18").Value = "工业富联"
$ws.Range("A19").Value = "天际股份"
$ws.Range("B19").Value = "新易盛"
$ws.Range("C19").Value = "方正科技"
$ws.Range("A20").Value = "万马股份"
$ws.Range("B20").Value = "卧龙电驱"
$ws.Range("C20").Value = "天赐材料"
$ws.Range("A21").Value = "天普股份"
$ws.Range("B21").Value = "幸福蓝海"
$ws.Range("C21").Value = "汇嘉时代"
